$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alt1")
Write-Host ("Dimension: " + $ws.Range("A1:R53").Address())
